{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change summary (from the target diff):\n//  1. The paragraph ending in \"\u5f9e\u5225\u7684\u7db2\u7ad9\u8b80\u53d6\u8a72\u7ad9\u8cc7\u6599\" loses the\n//     \"_GoBack\" bookmark that Word had auto-parked there.\n//  2. The empty paragraph right after the Google CDN link (just before\n//     the page break that starts the \"\u7db2\u7ad9\" heading) gets new content:\n//     a new bulleted (\"List Paragraph\" / numId 7) list item reading\n//     \"JSON (JavaScript Object Notation)\" \u2014 and the \"_GoBack\" bookmark\n//     now sits inside that new run, between \"Notation\" and \")\".\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- Locate the two target paragraphs by their (unique) text. ---\nlet corsParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (\n    paragraphs.items[i].text ===\n    \"\u7db2\u7ad9\u5982\u5141\u8a31CORS\uff0c\u5c31\u53ef\u4ee5\u900f\u904eAJAX\u5f9e\u5225\u7684\u7db2\u7ad9\u8b80\u53d6\u8a72\u7ad9\u8cc7\u6599\"\n  ) {\n    corsParaIndex = i;\n    break;\n  }\n}\n\n// The empty target paragraph (\"<w:p/>\" in the original) is the empty\n// paragraph that sits immediately before the page-break paragraph (whose\n// text starts with the form-feed / page-break character) that precedes\n// the \"\u7db2\u7ad9\" heading \u2014 i.e. the very last paragraph of the \"\u540d\u8a5e\"\n// (terms/glossary) section. Walking forward from the CORS paragraph and\n// matching on \"empty paragraph right before a page break\" uniquely and\n// robustly identifies it (there are several blank spacer paragraphs in\n// this section, and several page breaks elsewhere in the document, but\n// only one blank paragraph is directly followed by a page break here).\nlet jsonTargetIndex = -1;\nif (corsParaIndex !== -1) {\n  for (let i = corsParaIndex + 1; i < paragraphs.items.length - 1; i++) {\n    const t = paragraphs.items[i].text;\n    const nextText = paragraphs.items[i + 1].text;\n    if (t === \"\" && nextText.charCodeAt(0) === 12 /* page break */) {\n      jsonTargetIndex = i;\n      break;\n    }\n  }\n}\n\nif (corsParaIndex === -1 || jsonTargetIndex === -1) {\n  throw new Error(\n    \"Could not locate target paragraphs (cors=\" +\n      corsParaIndex +\n      \", json=\" +\n      jsonTargetIndex +\n      \")\"\n  );\n}\n\n// --- 1) Rewrite the CORS paragraph, dropping the _GoBack bookmark. ---\nconst corsPara = paragraphs.items[corsParaIndex];\nconst corsOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  \"<w:p>\" +\n  '<w:pPr><w:pStyle w:val=\"a9\"/><w:ind w:leftChars=\"0\" w:left=\"720\"/></w:pPr>' +\n  \"<w:r><w:t>\u7db2\u7ad9\u5982\u5141\u8a31</w:t></w:r>\" +\n  \"<w:r><w:t>CORS</w:t></w:r>\" +\n  \"<w:r><w:t>\uff0c\u5c31\u53ef\u4ee5\u900f\u904e</w:t></w:r>\" +\n  \"<w:r><w:t>AJAX</w:t></w:r>\" +\n  \"<w:r><w:t>\u5f9e\u5225\u7684\u7db2\u7ad9\u8b80\u53d6\u8a72\u7ad9\u8cc7\u6599</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\ncorsPara.insertOoxml(corsOoxml, Word.InsertLocation.replace);\n\n// --- 2) Turn the empty paragraph into the new \"JSON\" list item, with\n//        the _GoBack bookmark now living inside it. ---\nconst jsonPara = paragraphs.items[jsonTargetIndex];\nconst jsonOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  \"<w:p>\" +\n  '<w:pPr><w:pStyle w:val=\"a9\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"7\"/></w:numPr><w:ind w:leftChars=\"0\"/></w:pPr>' +\n  \"<w:r><w:t>JSON (JavaScript Object Notation</w:t></w:r>\" +\n  '<w:bookmarkStart w:id=\"1\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"1\"/>' +\n  \"<w:r><w:t>)</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\njsonPara.insertOoxml(jsonOoxml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is pre-seeded as $d below.\n#\n# Change summary (from the target diff):\n#  1. The paragraph ending in \"\u5f9e\u5225\u7684\u7db2\u7ad9\u8b80\u53d6\u8a72\u7ad9\u8cc7\u6599\" loses the\n#     \"_GoBack\" bookmark that Word had auto-parked there.\n#  2. The empty paragraph right after the Google CDN link (just before\n#     the page break that starts the \"\u7db2\u7ad9\" heading) gets new content:\n#     a new bulleted (\"List Paragraph\" / numId 7) list item reading\n#     \"JSON (JavaScript Object Notation)\" \u2014 and the \"_GoBack\" bookmark\n#     now sits inside that new run, between \"Notation\" and \")\".\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n# --- Locate the CORS paragraph by its (unique) text. ---\n$corsIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"\u7db2\u7ad9\u5982\u5141\u8a31CORS\uff0c\u5c31\u53ef\u4ee5\u900f\u904eAJAX\u5f9e\u5225\u7684\u7db2\u7ad9\u8b80\u53d6\u8a72\u7ad9\u8cc7\u6599\") {\n        $corsIdx = $i\n        break\n    }\n}\nif ($corsIdx -eq -1) {\n    throw \"Could not locate the CORS paragraph\"\n}\n\n# --- Locate the empty target paragraph (\"<w:p/>\" in the original OOXML).\n# It is the blank paragraph that sits immediately before the page-break\n# paragraph (whose text starts with the form-feed / page-break\n# character) which precedes the \"\u7db2\u7ad9\" heading \u2014 i.e. the very last\n# paragraph of the \"\u540d\u8a5e\" (terms/glossary) section. There are several\n# blank spacer paragraphs in this section and several page breaks\n# elsewhere in the document, but only one blank paragraph here is\n# directly followed by a page break, so this match is unique.\n$jsonIdx = -1\nfor ($i = $corsIdx + 1; $i -le $count - 1; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    $nextT = $d.Paragraphs.Item($i + 1).Range.Text\n    if ($t -eq \"\" -and $nextT.Length -gt 0 -and [int][char]$nextT[0] -eq 12) {\n        $jsonIdx = $i\n        break\n    }\n}\nif ($jsonIdx -eq -1) {\n    throw \"Could not locate the empty target paragraph\"\n}\n\n# --- 1) Rewrite the CORS paragraph, dropping the _GoBack bookmark. ---\n$corsOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"a9\"/><w:ind w:leftChars=\"0\" w:left=\"720\"/></w:pPr>' +\n  '<w:r><w:t>\u7db2\u7ad9\u5982\u5141\u8a31</w:t></w:r>' +\n  '<w:r><w:t>CORS</w:t></w:r>' +\n  '<w:r><w:t>\uff0c\u5c31\u53ef\u4ee5\u900f\u904e</w:t></w:r>' +\n  '<w:r><w:t>AJAX</w:t></w:r>' +\n  '<w:r><w:t>\u5f9e\u5225\u7684\u7db2\u7ad9\u8b80\u53d6\u8a72\u7ad9\u8cc7\u6599</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$corsRange = $d.Paragraphs.Item($corsIdx).Range\n$corsRange.InsertXML($corsOoxml)\n\n# --- 2) Turn the empty paragraph into the new \"JSON\" list item, with\n#        the _GoBack bookmark now living inside it. ---\n# Re-fetch the paragraph collection/index: inserting XML into the CORS\n# paragraph does not add or remove paragraphs, so $jsonIdx still points\n# at the right one.\n$jsonOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"a9\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"7\"/></w:numPr><w:ind w:leftChars=\"0\"/></w:pPr>' +\n  '<w:r><w:t>JSON (JavaScript Object Notation</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"1\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"1\"/>' +\n  '<w:r><w:t>)</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$jsonRange = $d.Paragraphs.Item($jsonIdx).Range\n$jsonRange.InsertXML($jsonOoxml)\n"}
